$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Populate the four new Salesforce "Account" test-scenario blocks
# (rows 2-14) that were added to the Test-Cases table.
# -----------------------------------------------------------------

# TestScenario_1 / TestCase_1 - New Account
$ws.Range("A2").Value = "TestScenario_1"
$ws.Range("B2").Value = "TestScenario_1.TestCase_1"
$ws.Range("C2").Value = "New Account"
$ws.Range("D2").Value = "User Needs to Login to Salesforce, from the browser with correct credentials"
$ws.Range("F2").Value = "Step 1"
$ws.Range("G2").Value = "Click on the Account tab, and click on New button"
$ws.Range("H2").Value = "User should be navigated to the New  Account Page"

$ws.Range("E3").Value = "Valid value for required field Name "
$ws.Range("F3").Value = "Step 2"
$ws.Range("G3").Value = "Input valid value in the  Name field."
$ws.Range("H3").Value = "User should be able to input value for the Name field."

$ws.Range("E4").Value = "Valid value for required field Annual Revenue "
$ws.Range("F4").Value = "Step 3"
$ws.Range("G4").Value = "Input valid value in the  Annual Revenue field."
$ws.Range("H4").Value = "User should be able to input value for the Annual Revenue field."

$ws.Range("F5").Value = "Step 4"
$ws.Range("G5").Value = "Click on Save button to save Account with fields"
$ws.Range("H5").Value = "User should be able to validate that a New Account is created"

# TestScenario_2 / TestCase_1 - View Account
$ws.Range("A6").Value = "TestScenario_2"
$ws.Range("B6").Value = "TestScenario_2.TestCase_1"
$ws.Range("C6").Value = "View Account"
$ws.Range("D6").Value = "User Needs to Login to Salesforce, from the browser with correct credentials"
$ws.Range("F6").Value = "Step 1"
$ws.Range("G6").Value = "Click on the Account tab,  and select a Account "
$ws.Range("H6").Value = "User should be navigated to the Account Page"

$ws.Range("F7").Value = "Step 2"
$ws.Range("G7").Value = "Click on Account name to View the Details"
$ws.Range("H7").Value = "User should be able to view the Account Details"

# TestScenario_3 / TestCase_1 - Edit Account
$ws.Range("A8").Value = "TestScenario_3"
$ws.Range("B8").Value = "TestScenario_3.TestCase_1"
$ws.Range("C8").Value = "Edit Account"
$ws.Range("D8").Value = "User Needs to Login to Salesforce, from the browser with correct credentials"
$ws.Range("F8").Value = "Step 1"
$ws.Range("G8").Value = "Click on the Account tab,  and click on existing Account to modify"
$ws.Range("H8").Value = "User is navigated to the Account Details page"

$ws.Range("E9").Value = "Valid value for required field Name "
$ws.Range("F9").Value = "Step 2"
$ws.Range("G9").Value = "Input valid value in the  Name field."
$ws.Range("H9").Value = "User should be able to input value for the Name field."

$ws.Range("E10").Value = "Valid value for required field Annual Revenue "
$ws.Range("F10").Value = "Step 3"
$ws.Range("G10").Value = "Input valid value in the  Annual Revenue field."
$ws.Range("H10").Value = "User should be able to input value for the Annual Revenue field."

$ws.Range("F11").Value = "Step 4"
$ws.Range("G11").Value = "Click on Save button to save Account with fields"
$ws.Range("H11").Value = "User should be able to validate that the Account is edited"

# TestScenario_4 / TestCase_1 - Delete Account
$ws.Range("A12").Value = "TestScenario_4"
$ws.Range("B12").Value = "TestScenario_4.TestCase_1"
$ws.Range("C12").Value = "Delete Account"
$ws.Range("D12").Value = "User Needs to Login to Salesforce, from the browser with correct credentials"
$ws.Range("F12").Value = "Step 1"
$ws.Range("G12").Value = "Click on the Account tab,  and select the existing  Account to delete"
$ws.Range("H12").Value = "User is navigated to the Account Details page"

$ws.Range("F13").Value = "Step 2"
$ws.Range("G13").Value = "Click on to the Delete to Delete the Account"
$ws.Range("H13").Value = "User should be able to validate that a pop-up is displayed asking for confirmation to delete the Account"

$ws.Range("F14").Value = "Step 3"
$ws.Range("G14").Value = "Click on Confirm / OK to delete the  Account"
$ws.Range("H14").Value = "User should be able to validate the Account is deleted"

# -----------------------------------------------------------------
# Grow Table1 so its range / AutoFilter cover the newly added rows
# -----------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J14"))

# -----------------------------------------------------------------
# Column widths as left by Excel after the new content was entered
# (values chosen so the engine's pixel-quantized ColumnWidth lands
# as close as possible to the recorded widths)
# -----------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 25.5
$ws.Columns.Item(3).ColumnWidth = 15.666666666666666
$ws.Columns.Item(4).ColumnWidth = 69.16666666666667
$ws.Columns.Item(5).ColumnWidth = 42.833333333333336
$ws.Columns.Item(7).ColumnWidth = 61.0
$ws.Columns.Item(8).ColumnWidth = 92.66666666666667
